$wb = $excel.ActiveWorkbook

# 1) Replace the "Ready for handoff" status text with "In Translation"
#    everywhere it appears across all worksheets.
foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    $rowCount = $used.Rows.Count
    $colCount = $used.Columns.Count
    for ($r = 1; $r -le $rowCount; $r++) {
        for ($c = 1; $c -le $colCount; $c++) {
            $cell = $ws.Cells.Item($r, $c)
            $txt = [string]$cell.Text
            if ($txt -eq "Ready for handoff") {
                $cell.Value2 = "In Translation"
            }
        }
    }
}

# 2) Narrow the "Status" column on each sheet (was ~17.22 chars, now ~13.41 chars).
$newWidth = 12.5

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E1:F1").ColumnWidth = $newWidth

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C1").ColumnWidth = $newWidth

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C1").ColumnWidth = $newWidth
